$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A85").Value = "2024-10-26 00:00:00"
$ws.Range("B85").Value = 73300
$ws.Range("C85").Value = 10273.44
$ws.Range("D85").Value = 9091.540000000001
$ws.Range("E85").Value = 7.1211
